$values = @(
    'basketball leg sleeve youth boys',
    'youth basketball leg sleeves boys',
    'basketball pants with knee pads',
    'mens compression pants',
    'compression knee pad',
    'basketball leggings',
    'knee compression pads',
    'compression pants men',
    'knee pads hex',
    'basketball knee pads youth boys',
    'mens capri pants',
    'volleyball knee',
    'workout hand pads',
    'black athletic leggings',
    'man pads',
    'compression tights',
    'tactical pants with knee pads',
    'mens pants',
    'xxl knee pads',
    'knee compression for men',
    'athletic leggings capri',
    'sliding pads for softball',
    'knee pads adult',
    'softball knee pads',
    'basketball knee pads adult',
    'youth basketball gear',
    'boys leggings',
    'compression knee pads',
    'black basketball leggings for men',
    'compression pants with knee pads boys',
    'hex knee pads basketball youth',
    'anti strip clothing for men',
    'hex protective knee pads',
    'men leggings',
    'weightlifting pads',
    'compression tights boys',
    'compression tights for boys',
    'bjj tights for men',
    'mens football pants',
    'advanced squat pad',
    'knee pads football',
    'knee protector pain',
    'knee support leggings',
    'compression tights for youth',
    'mens tall pants',
    'men tights',
    'wrestling gear',
    'little boys athletic leggings',
    'workout leggings for men pack',
    'baseball gear for boys',
    'black baseball pants youth xl',
    'wrestling knee pad',
    'volleyball knee pads xxl mens',
    'knee running',
    'knee pads for gym',
    'knee for running',
    'athletic compression leggings',
    'weightlifting pad',
    'foam knee pad',
    'compression tight',
    'knee support volleyball',
    'basketball equipment',
    'fitness squat pad',
    'knee swelling',
    'softball pants mens black',
    'sport pants for men',
    'basketball equipment pads',
    'knee length leggings',
    'volleyball knee pads girls youth',
    'reduce swelling after surgery',
    'soccer protection',
    'knee support for basketball men',
    'leg compression for men',
    'knee support for basketball',
    'soccer clothes',
    'lacrosse equipment',
    'capri athletic pants',
    'basketball knee pads youth pair',
    'mens leggings compression winter',
    'ski compression pants men',
    'leggings pants',
    'youth volleyball knee pads',
    'nike leggings basketball men',
    'underarmour tights',
    'yourh basketball knee pads',
    'compression leggings knee pads',
    'basketball tights with knee pads',
    'leggings for men with knee pads',
    'basketball knee pads men',
    'mens compression capri pants',
    'lacrosse leggings',
    'mens compression leggings capri',
    'basketball hip pads',
    'male compression tights',
    'compression tights for girls',
    'basketball knee pads',
    'mens compression pants capri',
    'compression mens pants',
    'basketball leggings for men',
    'youth capri compression tights'
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
